$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 13.45902438420715
$ws.Cells.Item(2, 3).Value = 8.528097289572054
$ws.Cells.Item(2, 5).Value = 19.79150204425236
$ws.Cells.Item(2, 6).Value = 40.13006075389559
$ws.Cells.Item(2, 7).Value = 32.33638601291391
$ws.Cells.Item(2, 8).Value = 14.82107250571386
$ws.Cells.Item(2, 10).Value = 8.014091543980664
$ws.Cells.Item(2, 13).Value = 18.94256346457001
$ws.Cells.Item(3, 2).Value = 12.88173193828261
$ws.Cells.Item(3, 3).Value = 8.038005303371792
$ws.Cells.Item(3, 5).Value = 19.77046571639309
$ws.Cells.Item(3, 6).Value = 40.04689846697832
$ws.Cells.Item(3, 7).Value = 32.19126801691338
$ws.Cells.Item(3, 8).Value = 14.87069728856357
$ws.Cells.Item(3, 10).Value = 8.044455614491246
$ws.Cells.Item(3, 13).Value = 18.7547492509905
$ws.Cells.Item(4, 2).Value = 12.51612410237143
$ws.Cells.Item(4, 3).Value = 7.719914472239958
$ws.Cells.Item(4, 5).Value = 19.76010073047676
$ws.Cells.Item(4, 6).Value = 40.00803581651642
$ws.Cells.Item(4, 7).Value = 32.11879389546263
$ws.Cells.Item(4, 8).Value = 14.90535328309774
$ws.Cells.Item(4, 10).Value = 8.06400557225265
$ws.Cells.Item(4, 13).Value = 18.64180006536371
$ws.Cells.Item(5, 2).Value = 12.36455813565574
$ws.Cells.Item(5, 3).Value = 7.586012827571999
$ws.Cells.Item(5, 5).Value = 19.75652184725252
$ws.Cells.Item(5, 6).Value = 39.99527183353349
$ws.Cells.Item(5, 7).Value = 32.09345080329171
$ws.Cells.Item(5, 8).Value = 14.92052288852982
$ws.Cells.Item(5, 10).Value = 8.072200961388942
$ws.Cells.Item(5, 13).Value = 18.59641252347251
$ws.Cells.Item(6, 2).Value = 12.33924218318894
$ws.Cells.Item(6, 3).Value = 7.563521657162964
$ws.Cells.Item(6, 5).Value = 19.75596662123224
$ws.Cells.Item(6, 6).Value = 39.99333806164358
$ws.Cells.Item(6, 7).Value = 32.08949580100947
$ws.Cells.Item(6, 8).Value = 14.92310485807988
$ws.Cells.Item(6, 10).Value = 8.073575629120782
$ws.Cells.Item(6, 13).Value = 18.58891589408523
$ws.Cells.Item(7, 2).Value = 12.51409014917219
$ws.Cells.Item(7, 3).Value = 7.71812588311631
$ws.Cells.Item(7, 5).Value = 19.76004984887407
$ws.Cells.Item(7, 6).Value = 40.00785123088571
$ws.Cells.Item(7, 7).Value = 32.1184351364141
$ws.Cells.Item(7, 8).Value = 14.90555363407589
$ws.Cells.Item(7, 10).Value = 8.064115171501808
$ws.Cells.Item(7, 13).Value = 18.64118530339398
$ws.Cells.Item(8, 2).Value = 13.2624142094714
$ws.Cells.Item(8, 3).Value = 8.362699091861574
$ws.Cells.Item(8, 5).Value = 19.78372127173961
$ws.Cells.Item(8, 6).Value = 40.09885861674835
$ws.Cells.Item(8, 7).Value = 32.28290631998611
$ws.Cells.Item(8, 8).Value = 14.83731095091015
$ws.Cells.Item(8, 10).Value = 8.024373413283611
$ws.Cells.Item(8, 13).Value = 18.87734026698517
$ws.Cells.Item(9, 2).Value = 14.63291465482997
$ws.Cells.Item(9, 3).Value = 9.489146979754043
$ws.Cells.Item(9, 5).Value = 19.85022786459036
$ws.Cells.Item(9, 6).Value = 40.37373227514463
$ws.Cells.Item(9, 7).Value = 32.73652664586432
$ws.Cells.Item(9, 8).Value = 14.73694713029713
$ws.Cells.Item(9, 10).Value = 7.953597445012477
$ws.Cells.Item(9, 13).Value = 19.35724548136507
$ws.Cells.Item(10, 2).Value = 15.57109545505918
$ws.Cells.Item(10, 3).Value = 10.23170607620581
$ws.Cells.Item(10, 5).Value = 19.9111110515796
$ws.Cells.Item(10, 6).Value = 40.63370847337481
$ws.Cells.Item(10, 7).Value = 33.14795514269333
$ws.Cells.Item(10, 8).Value = 14.68393009298633
$ws.Cells.Item(10, 10).Value = 7.905914260091892
$ws.Cells.Item(10, 13).Value = 19.71735858465026
$ws.Cells.Item(11, 2).Value = 15.98139369837426
$ws.Cells.Item(11, 3).Value = 10.55092244841084
$ws.Cells.Item(11, 5).Value = 19.94136478352881
$ws.Cells.Item(11, 6).Value = 40.76434623086038
$ws.Cells.Item(11, 7).Value = 33.35153933101523
$ws.Cells.Item(11, 8).Value = 14.66437492228718
$ws.Cells.Item(11, 10).Value = 7.885148961780028
$ws.Cells.Item(11, 13).Value = 19.88223319428863
$ws.Cells.Item(12, 2).Value = 16.13428623767779
$ws.Cells.Item(12, 3).Value = 10.66912455633096
$ws.Cells.Item(12, 5).Value = 19.95318376168888
$ws.Cells.Item(12, 6).Value = 40.8155695843891
$ws.Cells.Item(12, 7).Value = 33.43093467441599
$ws.Cells.Item(12, 8).Value = 14.65763086926088
$ws.Cells.Item(12, 10).Value = 7.877418084120746
$ws.Cells.Item(12, 13).Value = 19.94477001978994
$ws.Cells.Item(13, 2).Value = 16.10146986961446
$ws.Cells.Item(13, 3).Value = 10.64378680221687
$ws.Cells.Item(13, 5).Value = 19.9506222954819
$ws.Cells.Item(13, 6).Value = 40.80446016597873
$ws.Cells.Item(13, 7).Value = 33.41373408368785
$ws.Cells.Item(13, 8).Value = 14.65905384434605
$ws.Cells.Item(13, 10).Value = 7.879077185037298
$ws.Cells.Item(13, 13).Value = 19.93129782717796
$ws.Cells.Item(14, 2).Value = 15.99402249205547
$ws.Cells.Item(14, 3).Value = 10.56070070854159
$ws.Cells.Item(14, 5).Value = 19.94232990306255
$ws.Cells.Item(14, 6).Value = 40.76852541904957
$ws.Cells.Item(14, 7).Value = 33.35802545404711
$ws.Cells.Item(14, 8).Value = 14.66380680872111
$ws.Cells.Item(14, 10).Value = 7.884510286347634
$ws.Cells.Item(14, 13).Value = 19.88737636725412
$ws.Cells.Item(15, 2).Value = 15.92788217628129
$ws.Cells.Item(15, 3).Value = 10.50945915492762
$ws.Cells.Item(15, 5).Value = 19.93729762697467
$ws.Cells.Item(15, 6).Value = 40.74674191304329
$ws.Cells.Item(15, 7).Value = 33.32420027476685
$ws.Cells.Item(15, 8).Value = 14.66680436991402
$ws.Cells.Item(15, 10).Value = 7.887855452081273
$ws.Cells.Item(15, 13).Value = 19.86048507317971
$ws.Cells.Item(16, 2).Value = 15.54394016372692
$ws.Cells.Item(16, 3).Value = 10.21047004626632
$ws.Cells.Item(16, 5).Value = 19.90918492430625
$ws.Cells.Item(16, 6).Value = 40.62541785832121
$ws.Cells.Item(16, 7).Value = 33.13497601868824
$ws.Cells.Item(16, 8).Value = 14.68530031116468
$ws.Cells.Item(16, 10).Value = 7.907289897656226
$ws.Cells.Item(16, 13).Value = 19.7066006630892
$ws.Cells.Item(17, 2).Value = 15.30409971738655
$ws.Cells.Item(17, 3).Value = 10.02228577017968
$ws.Cells.Item(17, 5).Value = 19.89259013383685
$ws.Cells.Item(17, 6).Value = 40.55414141770537
$ws.Cells.Item(17, 7).Value = 33.0230596516728
$ws.Cells.Item(17, 8).Value = 14.69781917679295
$ws.Cells.Item(17, 10).Value = 7.919448998173388
$ws.Cells.Item(17, 13).Value = 19.61243291603187
$ws.Cells.Item(18, 2).Value = 15.16460350231801
$ws.Cells.Item(18, 3).Value = 9.912298698550263
$ws.Cells.Item(18, 5).Value = 19.88328632126216
$ws.Cells.Item(18, 6).Value = 40.51431147241187
$ws.Cells.Item(18, 7).Value = 32.96023818478715
$ws.Cells.Item(18, 8).Value = 14.70544886484531
$ws.Cells.Item(18, 10).Value = 7.92652979255722
$ws.Cells.Item(18, 13).Value = 19.55837259390896
$ws.Cells.Item(19, 2).Value = 15.11711046588717
$ws.Cells.Item(19, 3).Value = 9.874758934221383
$ws.Cells.Item(19, 5).Value = 19.88017776947012
$ws.Cells.Item(19, 6).Value = 40.50102678934895
$ws.Cells.Item(19, 7).Value = 32.93923579185577
$ws.Cells.Item(19, 8).Value = 14.70810569990482
$ws.Cells.Item(19, 10).Value = 7.928942225287885
$ws.Cells.Item(19, 13).Value = 19.54008781781929
$ws.Cells.Item(20, 2).Value = 15.32979201764861
$ws.Cells.Item(20, 3).Value = 10.04249932625947
$ws.Cells.Item(20, 5).Value = 19.89433176464799
$ws.Cells.Item(20, 6).Value = 40.56160837600403
$ws.Cells.Item(20, 7).Value = 33.03481334107488
$ws.Cells.Item(20, 8).Value = 14.69644206819011
$ws.Cells.Item(20, 10).Value = 7.918145621066077
$ws.Cells.Item(20, 13).Value = 19.6224469888136
$ws.Cells.Item(21, 2).Value = 16.02565043560946
$ws.Cells.Item(21, 3).Value = 10.58517783062498
$ws.Cells.Item(21, 5).Value = 19.94475578353804
$ws.Cells.Item(21, 6).Value = 40.77903295136633
$ws.Cells.Item(21, 7).Value = 33.37432644886332
$ws.Cells.Item(21, 8).Value = 14.66239276843638
$ws.Cells.Item(21, 10).Value = 7.882910862368185
$ws.Cells.Item(21, 13).Value = 19.90027477925615
$ws.Cells.Item(22, 2).Value = 16.46594874002066
$ws.Cells.Item(22, 3).Value = 10.92423818117149
$ws.Cells.Item(22, 5).Value = 19.979821358369
$ws.Cells.Item(22, 6).Value = 40.93134006312195
$ws.Cells.Item(22, 7).Value = 33.60960587919271
$ws.Cells.Item(22, 8).Value = 14.6439948905175
$ws.Cells.Item(22, 10).Value = 7.860654856871643
$ws.Cells.Item(22, 13).Value = 20.08242545861005
$ws.Cells.Item(23, 2).Value = 16.23231014306886
$ws.Cells.Item(23, 3).Value = 10.74470509879174
$ws.Cells.Item(23, 5).Value = 19.96091489085698
$ws.Cells.Item(23, 6).Value = 40.84912618926572
$ws.Cells.Item(23, 7).Value = 33.48282918379095
$ws.Cells.Item(23, 8).Value = 14.65345985093039
$ws.Cells.Item(23, 10).Value = 7.872462893055748
$ws.Cells.Item(23, 13).Value = 19.98517194714755
$ws.Cells.Item(24, 2).Value = 15.31818153122848
$ws.Cells.Item(24, 3).Value = 10.03336637565261
$ws.Cells.Item(24, 5).Value = 19.89354363552179
$ws.Cells.Item(24, 6).Value = 40.55822898687402
$ws.Cells.Item(24, 7).Value = 33.0294947567386
$ws.Cells.Item(24, 8).Value = 14.69706331238086
$ws.Cells.Item(24, 10).Value = 7.918734596318748
$ws.Cells.Item(24, 13).Value = 19.61791937985732
$ws.Cells.Item(25, 2).Value = 14.27362219676052
$ws.Cells.Item(25, 3).Value = 9.199287192800755
$ws.Cells.Item(25, 5).Value = 19.83010761845088
$ws.Cells.Item(25, 6).Value = 40.28911356802303
$ws.Cells.Item(25, 7).Value = 32.59992365036755
$ws.Cells.Item(25, 8).Value = 14.76048195556084
$ws.Cells.Item(25, 10).Value = 7.971982885597968
$ws.Cells.Item(25, 13).Value = 19.22589813612981
